# Adds blank paragraphs after the first four paragraphs and appends two new
# paragraphs ("Backoffice in LKP..." and "Read Times Techies...") each
# followed by a blank paragraph, per the target diff.

$WNS = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Insert-ParaAfter {
    param($d, [int]$paraIndex, [string]$innerXml)

    $anchor = $d.Paragraphs($paraIndex)
    $r = $anchor.Range
    $r.Collapse(0)            # wdCollapseEnd
    $r.InsertParagraphAfter()

    $newIndex = $paraIndex + 1
    $newPara = $d.Paragraphs($newIndex)
    $newPara.Range.InsertXML($innerXml)

    return $newIndex
}

$d = $word.ActiveDocument

# --- after "Use 5 Bars..." (paragraph 1) ---
$idx = Insert-ParaAfter $d 1 "<w:p $WNS/>"

# --- after "IFSC Code-11 Digits" ---
$idx = $idx + 1
$idx = Insert-ParaAfter $d $idx "<w:p $WNS/>"

# --- after "Acc num-After 220, 3 zeros" ---
$idx = $idx + 1
$idx = Insert-ParaAfter $d $idx "<w:p $WNS/>"

# --- after "To Book Auto ... copy location url" ---
$idx = $idx + 1
$idx = Insert-ParaAfter $d $idx "<w:p $WNS/>"

# --- new paragraph: Backoffice in LKP shows all the details related to Shares (2 runs) ---
$backofficeXml = "<w:p $WNS>" +
    '<w:r><w:t xml:space="preserve">Backoffice in LKP shows all the </w:t></w:r>' +
    '<w:r><w:t>details related to Shares</w:t></w:r>' +
    '</w:p>'
$idx = Insert-ParaAfter $d $idx $backofficeXml

# --- blank paragraph after Backoffice paragraph ---
$idx = Insert-ParaAfter $d $idx "<w:p $WNS/>"

# --- new paragraph: Read Times Techies on Thursday ---
$readTimesXml = "<w:p $WNS>" +
    '<w:r><w:t>Read Times Techies on Thursday</w:t></w:r>' +
    '</w:p>'
$idx = Insert-ParaAfter $d $idx $readTimesXml

# --- blank paragraph after Read Times Techies paragraph ---
$idx = Insert-ParaAfter $d $idx "<w:p $WNS/>"
